$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the "pi" block (rows 222-430) to "eta" ---
# (done first so that the new shared-string table assigns "eta" a lower index than "sig_mu")
for ($r = 222; $r -le 430; $r++) {
    $ws.Cells.Item($r, 1).Value = "eta"
}

# --- Rename the "V_mu" block (rows 211-221) to "sig_mu" and update its values ---
for ($r = 211; $r -le 221; $r++) {
    $ws.Cells.Item($r, 1).Value = "sig_mu"
}

$newValues = @(
    0.12210185334009473,
    0.11688998652357567,
    0.13201728054853046,
    0.11961147193854066,
    0.1261774832890005,
    0.14108522552013092,
    0.14902598494872799,
    0.12932332440483735,
    0.15581703349562331,
    0.0997142759384122,
    0.11483179428983986
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = 211 + $i
    $ws.Cells.Item($r, 4).Value = $newValues[$i]
}

# --- Update the view: scroll so row 193 is at top, select A211:A221 with active cell A211 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 193
$win.ScrollColumn = 1
$ws.Range("A211:A221").Select()
